$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("practica")

$ws.Range("D14").Value = 21
$ws.Range("D16").Value = 21
$ws.Range("D17").Value = 22
$ws.Range("D18").Value = 22
$ws.Range("D19").Value = 23
$ws.Range("D20").Value = 23
$ws.Range("D21").Value = 23
$ws.Range("D23").Value = 23
$ws.Range("D24").Value = 24
$ws.Range("D25").Value = 24
$ws.Range("D26").Value = 24
